$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 1989.3572
$ws.Range("I40").Value = 2037.75
$ws.Range("J40").Value = 1970
$ws.Range("K40").Value = 2037.75
$ws.Range("L40").Value = 1970
$ws.Range("M40").Value = -1862.75
$ws.Range("N40").Value = -2320

$ws.Range("H42").Value = 18
$ws.Range("I42").Value = 18
$ws.Range("J42").Value = 0
$ws.Range("K42").Value = 54
$ws.Range("L42").Value = 0
$ws.Range("M42").Value = 176
$ws.Range("N42").ClearContents()

$ws.Range("H116").Value = 4050
$ws.Range("J116").Value = 2500
$ws.Range("L116").Value = 2500
$ws.Range("N116").Value = -9384

$ws.Range("H132").Value = 20315.328
$ws.Range("I132").Value = 23386.105
$ws.Range("J132").Value = 2274.5
$ws.Range("K132").Value = 70158.315
$ws.Range("L132").Value = 6823.5
$ws.Range("M132").Value = -67628.315
$ws.Range("N132").Value = -11883.5

$ws.Range("H133").Value = 49477.5
$ws.Range("J133").Value = 49477.5
$ws.Range("L133").Value = 49477.5
$ws.Range("N133").Value = -59597.5

$ws.Range("H138").Value = 9692037
$ws.Range("I138").Value = 5105967.5
$ws.Range("J138").Value = 11118815
$ws.Range("K138").Value = 15317902.5
$ws.Range("L138").Value = 33356445
$ws.Range("M138").Value = -15312762.5
$ws.Range("N138").Value = -33366725

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 113215.89
$ws.Range("I2").Value = 145177.58
$ws.Range("J2").Value = 1350
$ws.Range("K2").Value = 145177.58
$ws.Range("L2").Value = 1350
$ws.Range("M2").Value = -145064.58
$ws.Range("N2").Value = -1576

$ws.Range("H4").Value = 1019.9167
$ws.Range("I4").Value = 1019.9167
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 1019.9167
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = -903.9167
$ws.Range("N4").ClearContents()

$ws.Range("H11").Value = 0
$ws.Range("I11").Value = 0
$ws.Range("K11").Value = 0
$ws.Range("M11").ClearContents()

$ws.Range("H32").Value = 64087
$ws.Range("I32").Value = 12403.7
$ws.Range("J32").Value = 322503.5
$ws.Range("K32").Value = 12403.7
$ws.Range("L32").Value = 322503.5
$ws.Range("M32").Value = -12116.7
$ws.Range("N32").Value = -323077.5

$ws.Range("H102").Value = 2379.158
$ws.Range("I102").Value = 2323.7646
$ws.Range("J102").Value = 2850
$ws.Range("K102").Value = 2323.7646
$ws.Range("L102").Value = 2850
$ws.Range("M102").Value = -701.7646
$ws.Range("N102").Value = -6094

$ws.Range("H116").Value = 113215.89
$ws.Range("I116").Value = 145177.58
$ws.Range("J116").Value = 1350
$ws.Range("K116").Value = 145177.58
$ws.Range("L116").Value = 1350
$ws.Range("M116").Value = -142883.58
$ws.Range("N116").Value = -5938

$ws.Range("H133").Value = 45632.5
$ws.Range("J133").Value = 45632.5
$ws.Range("L133").Value = 45632.5
$ws.Range("N133").Value = -50692.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 113215.89
$ws.Range("I3").Value = 145177.58
$ws.Range("J3").Value = 1350
$ws.Range("K3").Value = 145177.58
$ws.Range("L3").Value = 1350
$ws.Range("M3").Value = -145063.58
$ws.Range("N3").Value = -1578

$ws.Range("H105").Value = 3134.7368
$ws.Range("I105").Value = 2941.4814
$ws.Range("K105").Value = 2941.4814
$ws.Range("M105").Value = -1194.4814

$ws.Range("H123").Value = 29496
$ws.Range("J123").Value = 29496
$ws.Range("L123").Value = 29496
$ws.Range("N123").Value = -39296

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2050
$ws.Range("I16").Value = 1933.3334
$ws.Range("J16").Value = 2166.6667
$ws.Range("K16").Value = 1933.3334
$ws.Range("L16").Value = 2166.6667
$ws.Range("M16").Value = -1646.3334
$ws.Range("N16").Value = -2740.6667

$ws.Range("H107").Value = 512.34485
$ws.Range("I107").Value = 307.73914
$ws.Range("J107").Value = 1296.6666
$ws.Range("K107").Value = 307.73914
$ws.Range("L107").Value = 1296.6666
$ws.Range("M107").Value = 1612.26086
$ws.Range("N107").Value = -5136.6666

$ws.Range("H113").Value = 2050
$ws.Range("I113").Value = 1933.3334
$ws.Range("J113").Value = 2166.6667
$ws.Range("K113").Value = 1933.3334
$ws.Range("L113").Value = 2166.6667
$ws.Range("M113").Value = 236.6666
$ws.Range("N113").Value = -6506.6667

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 5953972
$ws.Range("I131").Value = 733.3333
$ws.Range("J131").Value = 6290948
$ws.Range("K131").Value = 2199.9999
$ws.Range("L131").Value = 18872844
$ws.Range("M131").Value = 2840.0001
$ws.Range("N131").Value = -18882924

$ws.Range("H137").Value = 11228788
$ws.Range("I137").Value = 16673510
$ws.Range("J137").Value = 339344.34
$ws.Range("K137").Value = 50020530
$ws.Range("L137").Value = 1018033.02
$ws.Range("M137").Value = -50015430
$ws.Range("N137").Value = -1028233.02

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 3291.3333
$ws.Range("I102").Value = 2111.1538
$ws.Range("J102").Value = 6359.8
$ws.Range("K102").Value = 2111.1538
$ws.Range("L102").Value = 6359.8
$ws.Range("M102").Value = -489.1538
$ws.Range("N102").Value = -9603.799999999999

$ws.Range("H126").Value = 2827.6316
$ws.Range("I126").Value = 2569.0667
$ws.Range("J126").Value = 2996.261
$ws.Range("K126").Value = 7707.2001
$ws.Range("L126").Value = 8988.782999999999
$ws.Range("M126").Value = -5237.2001
$ws.Range("N126").Value = -13928.783

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 11195.1
$ws.Range("I22").Value = 1393.875
$ws.Range("J22").Value = 50400
$ws.Range("K22").Value = 1393.875
$ws.Range("L22").Value = 50400
$ws.Range("M22").Value = -1098.875
$ws.Range("N22").Value = -50990

$ws.Range("H27").Value = 11195.1
$ws.Range("I27").Value = 1393.875
$ws.Range("J27").Value = 50400
$ws.Range("K27").Value = 1393.875
$ws.Range("L27").Value = 50400
$ws.Range("M27").Value = -1286.875
$ws.Range("N27").Value = -50614

$ws.Range("H55").Value = 845.3077
$ws.Range("I55").Value = 835.1818
$ws.Range("K55").Value = 835.1818
$ws.Range("M55").Value = -662.1818

$ws.Range("H61").Value = 1365.7646
$ws.Range("I61").Value = 1033.8
$ws.Range("J61").Value = 1840
$ws.Range("K61").Value = 1033.8
$ws.Range("L61").Value = 1840
$ws.Range("M61").Value = -831.8
$ws.Range("N61").Value = -2244

$ws.Range("H113").Value = 1365.7646
$ws.Range("I113").Value = 1033.8
$ws.Range("J113").Value = 1840
$ws.Range("K113").Value = 1033.8
$ws.Range("L113").Value = 1840
$ws.Range("M113").Value = 1136.2
$ws.Range("N113").Value = -6180

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 765.5769
$ws.Range("I113").Value = 643.56525
$ws.Range("J113").Value = 1701
$ws.Range("K113").Value = 1930.69575
$ws.Range("L113").Value = 5103
$ws.Range("M113").Value = 239.3042500000001
$ws.Range("N113").Value = -9443
